$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy column D (rows 2-4) values into new column E
$ws.Range("E2").Value = $ws.Range("D2").Value2
$ws.Range("E3").Value = $ws.Range("D3").Value2
$ws.Range("E4").Value = $ws.Range("D4").Value2

# Update the active selection as recorded in the saved file
$ws.Range("H10").Select()
